# Updates cryptos list values (price + 1h volume change) per the
# upstream data refresh, plus the Toncoin/Cardano row-order swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.399.76"
$ws.Range("E2").Value = "  -4.59%  "
$ws.Range("D3").Value = "2.375.32"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("D5").Value = "498.56"
$ws.Range("E5").Value = "  -7.09%  "
$ws.Range("D6").Value = "128.49"
$ws.Range("E6").Value = "  -4.62%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  -3.69%  "
$ws.Range("D9").Value = "2.396.86"
$ws.Range("E9").Value = "  -4.63%  "
$ws.Range("D10").Value = "0.0952"
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "4.66"
$ws.Range("E12").Value = "  -10.11%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "0.317"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("D14").Value = "2.799.14"
$ws.Range("E14").Value = "  -5.24%  "
$ws.Range("D15").Value = "56.261.43"
$ws.Range("E15").Value = "  -4.45%  "
$ws.Range("D16").Value = "21.43"
$ws.Range("E16").Value = "  -4.62%  "
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("D18").Value = "2.446.75"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.10"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.65%  "
$ws.Range("D20").Value = "309.28"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").Value = "4.02"
$ws.Range("E21").Value = "  -5.64%  "
$ws.Range("D22").Value = "6.23"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.00"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "2.495.14"
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("D27").Value = "0.372"
$ws.Range("E27").Value = "  -8.58%  "
$ws.Range("D28").Value = "0.151"
$ws.Range("E28").Value = "  -5.48%  "
$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").Value = "172.78"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").Value = "0.0₃0712"
$ws.Range("E31").Value = "  -6.59%  "
$ws.Range("D32").Value = "1.65"
$ws.Range("E32").Value = "  -5.09%  "
$ws.Range("D33").Value = "6.08"
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("E34").Value = "  -7.80%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "0.994"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "17.74"
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.20"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("D39").Value = "3.76"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("D40").Value = "35.83"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").Value = "0.788"
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "1.42"
$ws.Range("E42").Value = "  -6.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "130.00"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  -4.93%  "
$ws.Range("D45").Value = "4.85"
$ws.Range("E45").Value = "  -4.28%  "
$ws.Range("D46").Value = "0.572"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").Value = "252.02"
$ws.Range("E47").Value = "  -8.90%  "
$ws.Range("D48").Value = "0.0896"
$ws.Range("E48").Value = "  -4.92%  "
$ws.Range("D49").Value = "0.0481"
$ws.Range("E49").Value = "  -5.75%  "
$ws.Range("D50").Value = "16.75"
$ws.Range("E50").Value = "  -4.49%  "
$ws.Range("D51").Value = "0.0206"
$ws.Range("E51").Value = "  -6.19%  "
